# New 'Contractor' component (generalized 'CivicActions'); use secrender against templates
#
# Replace every occurrence of the "{{ organization_name }}" template token
# (the start of the "does not collect or maintain PII..." sentence) with
# the literal organization name "Example Org", leaving the rest of each
# paragraph (the paragraph's formatting/style plus the following " " and
# "control though..." runs) completely untouched.
#
# A plain Range.Text / Find.Execute replace would coalesce the whole
# paragraph into a single run - this engine, like real Word, re-flows a
# paragraph's runs whenever text inside it is deleted. InsertXML lets us
# swap in just the run content (keeping the existing <w:pPr> / style
# outside the replaced range, and copying the trailing runs verbatim), so
# only the placeholder text itself actually changes.

$d = $word.ActiveDocument

$oldToken = "{{ organization_name }}"
$newToken = "Example Org"
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs($i)
    $full = $para.Range
    $text = $full.Text

    if ($text.Length -ge $oldToken.Length -and $text.Substring(0, $oldToken.Length) -eq $oldToken) {
        # Drop the trailing paragraph-mark character Range.Text reports,
        # then apply the literal text substitution.
        $plain = $text -replace "[`r`a]+$", ""
        $newText = $newToken + $plain.Substring($oldToken.Length)

        # The known run layout for this sentence is:
        #   run 1: "<token ...> does not ... address this"
        #   run 2: " "
        #   run 3: "control though it may address it indirectly."
        # Rebuild the same three runs, only touching run 1's text, and
        # restrict the target Range to exclude both the paragraph mark
        # (End - 1) and the paragraph properties, so the original
        # <w:pPr>/style survives untouched.
        $marker = " control though it may address it indirectly."
        $markerIdx = $newText.IndexOf($marker)

        if ($markerIdx -ge 0) {
            $run1 = $newText.Substring(0, $markerIdx)
            $run3 = $newText.Substring($markerIdx + 1)
            $xml = "<w:p $wNs>" +
                   "<w:r><w:t xml:space=`"preserve`">$run1</w:t></w:r>" +
                   "<w:r><w:t xml:space=`"preserve`"> </w:t></w:r>" +
                   "<w:r><w:t xml:space=`"preserve`">$run3</w:t></w:r>" +
                   "</w:p>"
        } else {
            $xml = "<w:p $wNs><w:r><w:t xml:space=`"preserve`">$newText</w:t></w:r></w:p>"
        }

        $start = $full.Start
        $end = $full.End - 1
        $target = $d.Range($start, $end)
        $target.InsertXML($xml)
    }
}

$d.Save()
